$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Junio de 2020 a las 21:14"

# Row 4: Estados Unidos - updated stats
$ws.Range("B4").Value = 2104685
$ws.Range("C4").Value = 14984
$ws.Range("E4").Value = 1168826
$ws.Range("G4").Value = 450
$ws.Range("H4").Value = 116484

# Row 7: India - updated stats
$ws.Range("B7").Value = 309408
$ws.Range("C7").Value = 11125
$ws.Range("E7").Value = 146387

# Rows 76/77: Senegal <-> Uzbekistan swap places (with updated stats)
$ws.Range("A76").Value = "Uzbekistan"
$ws.Range("B76").Value = 4869
$ws.Range("C76").Value = 128
$ws.Range("D76").Value = 3700
$ws.Range("E76").Value = 1150
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 19

$ws.Range("A77").Value = "Senegal"
$ws.Range("B77").Value = 4851
$ws.Range("C77").Value = 92
$ws.Range("D77").Value = 3100
$ws.Range("E77").Value = 1695
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 56

# Rows 141/142: Jamaica <-> Yemen swap places (with updated stats)
$ws.Range("A141").Value = "Yemen"
$ws.Range("B141").Value = 632
$ws.Range("C141").Value = 41
$ws.Range("D141").Value = 28
$ws.Range("E141").Value = 465
$ws.Range("G141").Value = 3
$ws.Range("H141").Value = 139

$ws.Range("A142").Value = "Jamaica"
$ws.Range("B142").Value = 611
$ws.Range("C142").Value = 6
$ws.Range("D142").Value = 408
$ws.Range("E142").Value = 193
$ws.Range("H142").Value = 10

# Row 143: Togo - updated stats
$ws.Range("B143").Value = 525
$ws.Range("C143").Value = 1
$ws.Range("D143").Value = 279
$ws.Range("E143").Value = 233

# Rows 173/174: Camboya <-> Angola swap places (with updated stats)
$ws.Range("A173").Value = "Angola"
$ws.Range("B173").Value = 130
$ws.Range("C173").Value = 12
$ws.Range("D173").Value = 42
$ws.Range("E173").Value = 83
$ws.Range("H173").Value = 5

$ws.Range("A174").Value = "Camboya"
$ws.Range("B174").Value = 126
$ws.Range("D174").Value = 125
$ws.Range("E174").Value = 1
$ws.Range("H174").Value = 0

# Rows 206/207: Islas Malvinas <-> Groenlandia swap places (stats unchanged)
$ws.Range("A206").Value = "Groenlandia"
$ws.Range("A207").Value = "Islas Malvinas"

# Rows 208/209: Islas Turcas y Caicos <-> Santa Sede swap places (with updated stats)
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1

# Rows 210/211: Seychelles <-> Montserrat swap places (with updated stats)
$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# Rows 213/214: Papua Nueva Guinea <-> Islas Virgenes Britanicas swap places (with updated stats)
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
